$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPLKKEU029-001")

# Update the PREPARATION text in F2: change payment code and remove blank line
$ws.Range("F2").Value = "Username : 32193;`nPassword : bni1234;`nRole : 38 - Penyelia Teller;`nKode Pembayaran : DISK230200212"

# Update the KODE_PEMBAYARAN value in N2
$ws.Range("N2").Value = "DISK230200212"

# Update the Role value in I2: "Penyelia Teller" -> "Penyelia Settlement"
$ws.Range("I2").Value = "Penyelia Settlement"

# Update selection / view state to match new worksheet focus
$ws.Range("J2").Select()

# Columns D:L narrow slightly (Excel's "best fit" recompute once the role /
# payment-code text shrank) - set the new best-fit widths directly
$ws.Columns.Item(4).ColumnWidth = 15.5
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668
$ws.Columns.Item(6).ColumnWidth = 28.833333333333332
$ws.Columns.Item(7).ColumnWidth = 5.666666666666667
$ws.Columns.Item(8).ColumnWidth = 9
$ws.Columns.Item(9).ColumnWidth = 8.5
$ws.Columns.Item(10).ColumnWidth = 11.666666666666666
$ws.Columns.Item(11).ColumnWidth = 12.333333333333334
$ws.Columns.Item(12).ColumnWidth = 18.333333333333332

# Row 2 height shrinks now that the wrapped text in F2 has one less line
$ws.Rows.Item(2).RowHeight = 51
